$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Row 1: "100" -> "0M"
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"

# 2. Row 2: "0" -> "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"

# 3. Row 3: "17" -> "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"

# 4. Insert 10 new rows right after row 3 (before the current row 4),
#    preserving formatting from neighboring rows.
$newValues = @("11", "0.00003", "0.00004", "0.00004", "0.00000", "0.00003", "0.00004", "0.00004", "0.00041", "100.0")
$insertIndex = 4
foreach ($v in $newValues) {
    $refRow = $t.Rows.Item($insertIndex)
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $v
    $insertIndex = $insertIndex + 1
}

# After the insertion, the table has grown by 10 rows. The rows that used to
# hold the multi-run / tab-separated values are now at positions 44 and 45,
# and the trailing empty row is now at position 46.
$lastRowCount = $t.Rows.Count

# 5. Former row 34 (tab-separated "10 ... 100.0") -> single run "100"
$t.Rows.Item($lastRowCount - 2).Cells.Item(1).Range.Text = "100"

# 6. Former row 35 (tab-separated "1 ... 100.0") -> single run "0"
$t.Rows.Item($lastRowCount - 1).Cells.Item(1).Range.Text = "0"

# 7. Former row 36 (empty run) -> "17"
$t.Rows.Item($lastRowCount).Cells.Item(1).Range.Text = "17"
